$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 638
$ws.Range("F6").Value = 2822
$ws.Range("F8").Value = 56
$ws.Range("F10").Value = 582
$ws.Range("F11").Value = 23
$ws.Range("F12").Value = 325
$ws.Range("F14").Value = 5972
$ws.Range("F16").Value = 1048
$ws.Range("F17").Value = 10
$ws.Range("F18").Value = 235
$ws.Range("F19").Value = 173
$ws.Range("F21").Value = 540
$ws.Range("F22").Value = 32
$ws.Range("F24").Value = 102
$ws.Range("F25").Value = 1318
$ws.Range("F28").Value = 44
$ws.Range("F29").Value = 2060
$ws.Range("F30").Value = 181
$ws.Range("F31").Value = 354
$ws.Range("F33").Value = 3302

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 246
$ws.Range("F8").Value = 86
$ws.Range("F12").Value = 644
$ws.Range("G13").Value = 380
$ws.Range("F17").Value = 73
$ws.Range("F20").Value = 62
$ws.Range("F24").Value = 4054
$ws.Range("F25").Value = 9
$ws.Range("F29").Value = 216
$ws.Range("F36").Value = 20

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1803
$ws.Range("F5").Value = 2599
$ws.Range("F6").Value = 1142
$ws.Range("F8").Value = 1494
$ws.Range("F12").Value = 644

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1803
$ws.Range("F4").Value = 2599
$ws.Range("F6").Value = 1142
$ws.Range("F7").Value = 1494
$ws.Range("F11").Value = 638
$ws.Range("F12").Value = 2822
$ws.Range("F13").Value = 56
$ws.Range("F14").Value = 644
$ws.Range("F15").Value = 582
$ws.Range("F16").Value = 86
$ws.Range("F17").Value = 23
$ws.Range("F18").Value = 325
$ws.Range("F20").Value = 5972
$ws.Range("F23").Value = 1048
$ws.Range("F24").Value = 235
$ws.Range("F25").Value = 173
$ws.Range("F27").Value = 540
$ws.Range("F29").Value = 73
$ws.Range("F31").Value = 62
$ws.Range("F36").Value = 9
$ws.Range("F40").Value = 216
$ws.Range("F41").Value = 44
$ws.Range("F44").Value = 2060
$ws.Range("F47").Value = 181
$ws.Range("F48").Value = 354
$ws.Range("F49").Value = 20
$ws.Range("F50").Value = 3302
